$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.869.27"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "2.353.24"
$ws.Range("E3").Value = "  -2.61%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'320.58"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").Value = "'105.18"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.637"
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("D9").Value = "'0.617"
$ws.Range("E9").Value = "  -6.60%  "
$ws.Range("D10").Value = "'41.14"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("D11").Value = "'0.0924"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("D12").Value = "'8.48"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").Value = "'0.997"
$ws.Range("E13").Value = "  -4.80%  "
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").Value = "'15.96"
$ws.Range("E15").Value = "  -8.85%  "
$ws.Range("D16").Value = "2.712.33"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "2.332.54"
$ws.Range("E17").Value = "  -8.73%  "
$ws.Range("D18").Value = "42.854.97"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").Value = "'7.72"
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("E20").Value = "  -4.17%  "
$ws.Range("D21").Value = "'77.31"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").Value = "'3.62"
$ws.Range("E22").Value = "  +3.40%  "
$ws.Range("D23").Value = "'260.36"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = "  -5.27%  "
$ws.Range("D25").Value = "'9.54"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'11.41"
$ws.Range("E27").Value = "  -4.92%  "
$ws.Range("D28").Value = "'23.26"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").Value = "'174.57"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("D31").Value = "'36.27"
$ws.Range("E31").Value = "  -4.62%  "
$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "'3.00"
$ws.Range("E32").Value = "  -6.82%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0892"
$ws.Range("E33").Value = "  -4.40%  "
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("D36").Value = "'0.114"
$ws.Range("E36").Value = "  +6.60%  "
$ws.Range("D37").Value = "'4.62"
$ws.Range("E37").Value = "  -5.72%  "
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").Value = "'3.80"
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("D40").Value = "'2.68"
$ws.Range("E40").Value = "  -7.54%  "
$ws.Range("D41").Value = "'71.69"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("E42").Value = "  -10.27%  "
$ws.Range("D43").Value = "'0.232"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "'115.09"
$ws.Range("E45").Value = "  -9.49%  "
$ws.Range("D46").Value = "'11.83"
$ws.Range("E46").Value = "  -6.51%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'5.52"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'86.33"
$ws.Range("E48").Value = "  +4.16%  "
$ws.Range("D49").Value = "'9.22"
$ws.Range("E49").Value = "  -5.05%  "
$ws.Range("D50").Value = "'73.63"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("E51").Value = "  -1.17%  "
